$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each written cell to retain its original General/Normal style while
# storing the new value as literal text (matches the source data which used
# inline strings like "87.706.68" / "1.00" / "  +1.63%  " that must not be
# reinterpreted as numbers).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '88.438.67'
Set-TextValue 'E2' '  +2.62%  '
Set-TextValue 'D3' '3.187.57'
Set-TextValue 'E3' '  -2.43%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  -0.43%  '
Set-TextValue 'D5' '209.55'
Set-TextValue 'E5' '  +0.04%  '
Set-TextValue 'D6' '615.26'
Set-TextValue 'E6' '  -1.40%  '
Set-TextValue 'D7' '0.398'
Set-TextValue 'E7' '  +7.88%  '
Set-TextValue 'D8' '0.684'
Set-TextValue 'E8' '  +5.24%  '
Set-TextValue 'D9' '0.997'
Set-TextValue 'E9' '  -0.31%  '
Set-TextValue 'D10' '3.175.78'
Set-TextValue 'E10' '  -2.71%  '
Set-TextValue 'D11' '0.546'
Set-TextValue 'E11' '  -5.19%  '
Set-TextValue 'D12' '0.177'
Set-TextValue 'E12' '  -1.12%  '
Set-TextValue 'D13' '0.0000253'
Set-TextValue 'E13' '  -1.11%  '
Set-TextValue 'D14' '5.33'
Set-TextValue 'E14' '  +1.21%  '
Set-TextValue 'D15' '3.778.49'
Set-TextValue 'E15' '  -2.30%  '
Set-TextValue 'B16' 'Avalanche'
Set-TextValue 'C16' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D16' '32.80'
Set-TextValue 'E16' '  -3.00%  '
Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '88.035.96'
Set-TextValue 'E17' '  +2.35%  '
Set-TextValue 'D18' '3.165.03'
Set-TextValue 'E18' '  -2.81%  '
Set-TextValue 'D19' '3.20'
Set-TextValue 'E19' '  +7.83%  '
Set-TextValue 'D20' '13.62'
Set-TextValue 'E20' '  -2.94%  '
Set-TextValue 'D21' '416.27'
Set-TextValue 'E21' '  -2.89%  '
Set-TextValue 'D22' '8.55'
Set-TextValue 'E22' '  -4.71%  '
Set-TextValue 'D23' '5.16'
Set-TextValue 'E23' '  -2.36%  '
Set-TextValue 'B24' 'PEPE'
Set-TextValue 'C24' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D24' '0.0000163'
Set-TextValue 'E24' '  +26.75%  '
Set-TextValue 'B25' 'NEARProtocol'
Set-TextValue 'C25' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D25' '5.31'
Set-TextValue 'E25' '  +3.18%  '
Set-TextValue 'B26' 'Aptos'
Set-TextValue 'C26' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D26' '12.37'
Set-TextValue 'E26' '  +0.99%  '
Set-TextValue 'B27' 'WrappedeETH'
Set-TextValue 'C27' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D27' '3.333.63'
Set-TextValue 'E27' '  -2.96%  '
Set-TextValue 'D28' '73.94'
Set-TextValue 'E28' '  -2.78%  '
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  -0.19%  '
Set-TextValue 'D30' '0.167'
Set-TextValue 'E30' '  -3.68%  '
Set-TextValue 'E31' '  +0.20%  '
Set-TextValue 'D32' '553.57'
Set-TextValue 'E32' '  +1.70%  '
Set-TextValue 'D33' '8.38'
Set-TextValue 'E33' '  -4.79%  '
Set-TextValue 'D34' '1.33'
Set-TextValue 'E34' '  -6.03%  '
Set-TextValue 'D35' '6.96'
Set-TextValue 'E35' '  +2.28%  '
Set-TextValue 'D36' '1.87'
Set-TextValue 'E36' '  -3.86%  '
Set-TextValue 'D37' '0.132'
Set-TextValue 'E37' '  -3.53%  '
Set-TextValue 'D38' '22.06'
Set-TextValue 'E38' '  -1.55%  '
Set-TextValue 'B39' 'WhiteBITCoin'
Set-TextValue 'C39' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D39' '21.81'
Set-TextValue 'E39' '  +0.78%  '
Set-TextValue 'B40' 'dogwifhat'
Set-TextValue 'C40' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D40' '3.20'
Set-TextValue 'E40' '  +10.16%  '
Set-TextValue 'B41' 'FirstDigitalUSD'
Set-TextValue 'C41' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D41' '0.995'
Set-TextValue 'E41' '  -0.49%  '
Set-TextValue 'E42' '  +0.05%  '
Set-TextValue 'D43' '1.94'
Set-TextValue 'E43' '  -2.31%  '
Set-TextValue 'D44' '0.379'
Set-TextValue 'E44' '  -3.30%  '
Set-TextValue 'D45' '150.29'
Set-TextValue 'E45' '  -4.14%  '
Set-TextValue 'D46' '175.63'
Set-TextValue 'E46' '  -0.44%  '
Set-TextValue 'D47' '43.16'
Set-TextValue 'E47' '  -1.88%  '
Set-TextValue 'D48' '0.126'
Set-TextValue 'E48' '  +5.81%  '
Set-TextValue 'D49' '1.25'
Set-TextValue 'E49' '  -4.94%  '
Set-TextValue 'D50' '24.43'
Set-TextValue 'E50' '  +1.53%  '
Set-TextValue 'D51' '4.01'
Set-TextValue 'E51' '  -5.94%  '
